$d = $word.ActiveDocument

# --- Change 1: "Tecnologias utilizadas" paragraph -----------------------
# Before: " Java" + " y" + " MySQL (JDBC)" + ";" + " " + "Fronte" + "n" + "d" + ...
# After:  " Java" + " 17 LTS" + " y" + " MySQL" + " 8.2, y" + " " + "Fronte" + "n" + "d" + ...
# Do this as a single Find/Replace across the affected span so the
# surrounding (unchanged) runs on both sides are left alone.
$r1 = $d.Content
$r1.Find.Execute(
    " Java y MySQL (JDBC);",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    " Java 17 LTS y MySQL 8.2, y",
    2)

# --- Change 2: "Sprint 3 - ..." heading paragraph ------------------------
# Merge the three bold runs that make up the Sprint 3 heading into a single
# run with the same text, exactly like the target document does.
$r2 = $d.Content
$r2.Find.Execute(
    "Sprint 3 - Creación de la base de datos y la estructura arquitectónica del proyecto (1 Semana) 15 – 21 enero 2024",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Sprint 3 - Creación de la base de datos y la estructura arquitectónica del proyecto (1 Semana) 15 – 21 enero 2024",
    2)
